$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Clec11a/Itga11 -> FAPs (counts updated: pair count 1 -> 3)
$row2 = @("FAPs","Clec11a","Itga11","FAPs",3,1,14.768619,44.305857,0.9736679609684162,0.9736679609684162,3,1,24.759128,74.277384,0.9895671066967037,0.9895671066967037,365.658128204232,3290.923153838088,0.9635097870187946,0.9635097870187946)

# Row 3: FAPs -> Clec11a/Itga11 -> sCs (new target cluster, re-derived receptor stats)
$row3 = @("FAPs","Clec11a","Itga11","sCs",3,1,14.768619,44.305857,0.9736679609684162,0.9736679609684162,2,0.6666666666666666,0.2610326666666667,0.7830980000000001,0.0104328933032964,0.0104328933032964,3.855092000554001,34.69582800498601,0.01015817394962165,0.01015817394962165)

# Row 4 (new): sCs -> Clec11a/Itga11 -> FAPs
$row4 = @("sCs","Clec11a","Itga11","FAPs",3,1,0.399405,1.198215,0.02633203903158381,0.02633203903158381,3,1,24.759128,74.277384,0.9895671066967037,0.9895671066967037,9.88891951884,89.00027566956,0.02605731967790906,0.02605731967790906)

# Row 5 (new): sCs -> Clec11a/Itga11 -> sCs
$row5 = @("sCs","Clec11a","Itga11","sCs",3,1,0.399405,1.198215,0.02633203903158381,0.02633203903158381,2,0.6666666666666666,0.2610326666666667,0.7830980000000001,0.0104328933032964,0.0104328933032964,0.10425775223,0.9383197700700001,0.0002747193536747502,0.0002747193536747502)

$rows = @($row2, $row3, $row4, $row5)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $vals = $rows[$i]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
